$wb = $excel.ActiveWorkbook

# --- Keywords sheet: add new "Proceed_on_fail" column (F) ---
$ws = $wb.Worksheets.Item("Keywords")

# Copy formatting (style/fill) from the existing header cell E1 into new header F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "Proceed_on_fail"
$ws.Range("F3").Value = "N"
$ws.Range("F4").Value = "N"
$ws.Range("F5").Value = "N"
$ws.Range("F6").Value = "N"
$ws.Range("F7").Value = "N"
$ws.Range("F8").Value = "N"

$ws.Range("F2").Select()

# --- Make "Keywords" the active (selected/visible) tab ---
$ws.Select()

# --- "Test Cases" sheet keeps its own selection ---
$wsTC = $wb.Worksheets.Item("Test Cases")
$wsTC.Range("B2").Select()

# --- Switch back to the Keywords sheet so it stays the active tab on save ---
$ws.Select()
